$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 5.633000000000001
$ws.Range("C3").Value = -12.787
$ws.Range("A12").Value = -21.66
$ws.Range("B14").Value = 5.929
$ws.Range("C20").Value = -12.46
$ws.Range("C25").Value = -12.738
$ws.Range("B26").Value = 6.419
$ws.Range("A27").Value = -21.217
$ws.Range("C30").Value = -13.734
$ws.Range("B31").Value = 6.205
$ws.Range("A32").Value = -21.471
$ws.Range("B35").Value = 7.342000000000001
$ws.Range("A36").Value = -20.934
$ws.Range("B37").Value = 7.641
$ws.Range("A38").Value = -20.363
$ws.Range("C44").Value = -12.646
$ws.Range("B45").Value = 5.747
$ws.Range("A46").Value = -21.591
$ws.Range("C47").Value = -12.392
$ws.Range("B52").Value = 5.4
$ws.Range("A54").Value = -21.585
$ws.Range("A55").Value = -21.825
$ws.Range("A56").Value = -21.832
$ws.Range("B57").Value = 6.090000000000001
$ws.Range("C58").Value = -12.802
$ws.Range("A67").Value = -21.588
$ws.Range("A69").Value = -21.721
$ws.Range("A72").Value = -21.567
$ws.Range("C78").Value = -13.032
$ws.Range("B81").Value = 6.693
$ws.Range("A83").Value = -20.489
$ws.Range("B83").Value = 7.256
$ws.Range("C84").Value = -13.509
$ws.Range("A86").Value = -22.179
$ws.Range("C89").Value = -11.275
$ws.Range("A91").Value = -21.747
$ws.Range("C91").Value = -11.21
$ws.Range("C92").Value = -11.539
$ws.Range("A93").Value = -21.593
$ws.Range("C96").Value = -13.261
$ws.Range("A99").Value = -20.683
$ws.Range("B100").Value = 6.218999999999999
$ws.Range("B102").Value = 6.667
$ws.Range("C102").Value = -12.798
